$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New daily COVID data rows (2020-05-25, 2020-05-26, 2020-05-27).
$newData = @(
    @(43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
    @(43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0)
)

# Insert the new rows right after the current last data row (75), copying
# that row's formatting (styles) down so the table keeps a consistent look.
$lastRow = 75
for ($i = 0; $i -lt $newData.Count; $i++) {
    $newRowIndex = $lastRow + $i + 1
    $null = $ws.Rows($lastRow).Copy()
    $null = $ws.Rows($newRowIndex).Insert()
}

# Grow the table (ListObject) - and therefore its AutoFilter range - to
# cover the newly inserted rows.
$newLastRow = $lastRow + $newData.Count
$null = $lo.Resize($ws.Range("A1:J" + $newLastRow))

# Fill in the values for the new rows.
$r = $lastRow + 1
foreach ($rowVals in $newData) {
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
    $r = $r + 1
}

# Match the selection Excel leaves behind after entering the last new row.
$null = $ws.Range("A" + $newLastRow + ":J" + $newLastRow).Select()
